$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add a new "f" / "lek 18" column (F) for a handful of rows -------------
# Rows 2, 4, 7 and 8 gain a new F cell containing "f".
$ws.Range("F2").Value = "f"
$ws.Range("F4").Value = "f"
$ws.Range("F7").Value = "f"
$ws.Range("F8").Value = "f"

# --- Collapse the Lek11 / Lek12 / Lek13 variants down to "Lek10" -----------
$ws.Range("A12").Value = "Lek10"
$ws.Range("A13").Value = "Lek10"
$ws.Range("A24").Value = "Lek10"
$ws.Range("A25").Value = "Lek10"

# --- Replace the scratch/test row 26 with a single "lek 18" label ----------
$ws.Range("B26:E26").ClearContents()
$ws.Range("A26").Value = "lek 18"

# --- Window / selection bookkeeping (matches the saved view state) ---------
$ws.Range("E27").Select()
